$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4450399954297666
$ws.Range("C2").Value = 0.2617755575175842
$ws.Range("E2").Value = 0.1473957609363019
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.5915530940608633
$ws.Range("H2").Value = 0.7338712398171623
$ws.Range("K2").Value = 0.2191570370909233
$ws.Range("L2").Value = 0.1894865821604412
$ws.Range("M2").Value = 0.1380051491995147
$ws.Range("N2").Value = 1.735146594377166
$ws.Range("O2").Value = 2.631368660826297

$ws.Range("B3").Value = 0.4125773520567861
$ws.Range("C3").Value = 0.2629422262831156
$ws.Range("E3").Value = 0.1479537794625685
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.5965272105873609
$ws.Range("H3").Value = 0.7393926312946917
$ws.Range("K3").Value = 0.1915746033014472
$ws.Range("L3").Value = 0.1868115878728389
$ws.Range("M3").Value = 0.1316839357053077
$ws.Range("N3").Value = 1.750151953579714
$ws.Range("O3").Value = 2.653261325625095

$ws.Range("B4").Value = 0.3927459779189064
$ws.Range("C4").Value = 0.2637070686545115
$ws.Range("E4").Value = 0.1483565256792083
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.5999130226853993
$ws.Range("H4").Value = 0.7430431519180374
$ws.Range("K4").Value = 0.1745924537895007
$ws.Range("L4").Value = 0.1852572552978131
$ws.Range("M4").Value = 0.1278529675646638
$ws.Range("N4").Value = 1.759886326140663
$ws.Range("O4").Value = 2.667944423447992

$ws.Range("B5").Value = 0.3846905352847898
$ws.Range("C5").Value = 0.2640309765569313
$ws.Range("E5").Value = 0.1485357946542116
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.6013761829536577
$ws.Range("H5").Value = 0.7445963109172311
$ws.Range("K5").Value = 0.1676608909734796
$ws.Range("L5").Value = 0.1846460965908463
$ws.Range("M5").Value = 0.1263045974731405
$ws.Range("N5").Value = 1.763984274442475
$ws.Range("O5").Value = 2.674240107995416

$ws.Range("B6").Value = 0.3833545248660073
$ws.Range("C6").Value = 0.2640855007264058
$ws.Range("E6").Value = 0.1485664777959101
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.6016241779058902
$ws.Range("H6").Value = 0.7448581724806544
$ws.Range("K6").Value = 0.1665092478193628
$ws.Range("L6").Value = 0.1845459601526969
$ws.Range("M6").Value = 0.1260482677274695
$ws.Range("N6").Value = 1.764672657399778
$ws.Range("O6").Value = 2.675304362775833

$ws.Range("B7").Value = 0.3926372332530548
$ws.Range("C7").Value = 0.2637113874356736
$ws.Range("E7").Value = 0.1483588819967387
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.5999324176217655
$ws.Range("H7").Value = 0.7430638328899022
$ws.Range("K7").Value = 0.1744990170251413
$ws.Range("L7").Value = 0.1852489228182534
$ws.Range("M7").Value = 0.1278320337750962
$ws.Range("N7").Value = 1.759941061504396
$ws.Range("O7").Value = 2.668028064932514

$ws.Range("B8").Value = 0.4338263310911259
$ws.Range("C8").Value = 0.2621677773037518
$ws.Range("E8").Value = 0.1475757048354023
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.5931993339262789
$ws.Range("H8").Value = 0.7357210201896365
$ws.Range("K8").Value = 0.2096565256678531
$ws.Range("L8").Value = 0.1885460034658664
$ws.Range("M8").Value = 0.1358152384374165
$ws.Range("N8").Value = 1.740212363779648
$ws.Range("O8").Value = 2.6386597993028

$ws.Range("B9").Value = 0.5153734635564433
$ws.Range("C9").Value = 0.2595241822427674
$ws.Range("E9").Value = 0.1465157295003401
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.5826273064104086
$ws.Range("H9").Value = 0.7233844621923993
$ws.Range("K9").Value = 0.2782142747421403
$ws.Range("L9").Value = 0.1957075893008167
$ws.Range("M9").Value = 0.1518641099739924
$ws.Range("N9").Value = 1.705653739300804
$ws.Range("O9").Value = 2.590907838312006

$ws.Range("B10").Value = 0.5757319247110786
$ws.Range("C10").Value = 0.2578137131661649
$ws.Range("E10").Value = 0.1460255376490807
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.5764642414502603
$ws.Range("H10").Value = 0.7155738835081635
$ws.Range("K10").Value = 0.3283292669257207
$ws.Range("L10").Value = 0.2013899631244982
$ws.Range("M10").Value = 0.1638896720886436
$ws.Range("N10").Value = 1.682773679330793
$ws.Range("O10").Value = 2.561813725054179

$ws.Range("B11").Value = 0.6032819815232813
$ws.Range("C11").Value = 0.2570854904100877
$ws.Range("E11").Value = 0.1458648844463823
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.5740088548473068
$ws.Range("H11").Value = 0.712291826649988
$ws.Range("K11").Value = 0.3510688649908218
$ws.Range("L11").Value = 0.2040656742335472
$ws.Range("M11").Value = 0.169410199228814
$ws.Range("N11").Value = 1.672908493155649
$ws.Range("O11").Value = 2.549876740409573

$ws.Range("B12").Value = 0.6137272110837557
$ws.Range("C12").Value = 0.2568168716567278
$ws.Range("E12").Value = 0.145812985856498
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.5731291374903478
$ws.Range("H12").Value = 0.7110878996335757
$ws.Range("K12").Value = 0.3596709979224784
$ws.Range("L12").Value = 0.2050918727856441
$ws.Range("M12").Value = 0.171507751461867
$ws.Range("N12").Value = 1.669250804221054
$ws.Range("O12").Value = 2.545543036101563

$ws.Range("B13").Value = 0.6114770939618097
$ws.Range("C13").Value = 0.2568744062676949
$ws.Range("E13").Value = 0.1458237660845576
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.5733163726938741
$ws.Range("H13").Value = 0.7113454568978526
$ws.Range("K13").Value = 0.357818777054149
$ws.Range("L13").Value = 0.2048702873546944
$ws.Range("M13").Value = 0.1710556949223019
$ws.Range("N13").Value = 1.670035082761519
$ws.Range("O13").Value = 2.5464680804088

$ws.Range("B14").Value = 0.6041410679575847
$ws.Range("C14").Value = 0.2570632479835808
$ws.Range("E14").Value = 0.1458604357549156
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.5739354761504103
$ws.Range("H14").Value = 0.7121919992143759
$ws.Range("K14").Value = 0.3517767482021554
$ws.Range("L14").Value = 0.2041498408756439
$ws.Range("M14").Value = 0.1695826256958384
$ws.Range("N14").Value = 1.672606009064367
$ws.Range("O14").Value = 2.549516465261974

$ws.Range("B15").Value = 0.5996491628164335
$ws.Range("C15").Value = 0.2571798484133971
$ws.Range("E15").Value = 0.1458840600909461
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.5743212176405379
$ws.Range("H15").Value = 0.7127155973050634
$ws.Range("K15").Value = 0.3480746617604211
$ws.Range("L15").Value = 0.2037102323948261
$ws.Range("M15").Value = 0.1686812418039239
$ws.Range("N15").Value = 1.674190937724802
$ws.Range("O15").Value = 2.551407983643102

$ws.Range("B16").Value = 0.5739332683750433
$ws.Range("C16").Value = 0.2578623054283931
$ws.Range("E16").Value = 0.1460372888164869
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.5766317147565019
$ws.Range("H16").Value = 0.7157938224493847
$ws.Range("K16").Value = 0.326841969528374
$ws.Range("L16").Value = 0.2012169181429897
$ws.Range("M16").Value = 0.1635298869939277
$ws.Range("N16").Value = 1.683429316077063
$ws.Range("O16").Value = 2.562619948056593

$ws.Range("B17").Value = 0.5581806367636659
$ws.Range("C17").Value = 0.2582937249867214
$ws.Range("E17").Value = 0.1461472367745777
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.5781383255055843
$ws.Range("H17").Value = 0.717751588046994
$ws.Range("K17").Value = 0.3138011817716801
$ws.Range("L17").Value = 0.1997105367191665
$ws.Range("M17").Value = 0.1603824091767763
$ws.Range("N17").Value = 1.689235818723123
$ws.Range("O17").Value = 2.569830529800484

$ws.Range("B18").Value = 0.5491289064633804
$ws.Range("C18").Value = 0.2585465626030157
$ws.Range("E18").Value = 0.1462163450381091
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.5790376639294621
$ws.Range("H18").Value = 0.7189031554695262
$ws.Range("K18").Value = 0.306295041139748
$ws.Range("L18").Value = 0.1988526550116347
$ws.Range("M18").Value = 0.1585767818601411
$ws.Range("N18").Value = 1.692626684405088
$ws.Range("O18").Value = 2.574100051837007

$ws.Range("B19").Value = 0.5460656733731923
$ws.Range("C19").Value = 0.258632976590274
$ws.Range("E19").Value = 0.1462407527484686
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.5793477934165594
$ws.Range("H19").Value = 0.7192974398757244
$ws.Range("K19").Value = 0.3037526752183055
$ws.Range("L19").Value = 0.1985636619871514
$ws.Range("M19").Value = 0.1579662423911188
$ws.Range("N19").Value = 1.693783554535088
$ws.Range("O19").Value = 2.575566626922011

$ws.Range("B20").Value = 0.5598566290896088
$ws.Range("C20").Value = 0.258247313779453
$ws.Range("E20").Value = 0.1461349253903172
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.5779745519747124
$ws.Range("H20").Value = 0.7175405405147686
$ws.Range("K20").Value = 0.3151899610963653
$ws.Range("L20").Value = 0.1998700094352728
$ws.Range("M20").Value = 0.1607169762631102
$ws.Range("N20").Value = 1.688612416351543
$ws.Range("O20").Value = 2.569050305900404

$ws.Range("B21").Value = 0.6062954997366603
$ws.Range("C21").Value = 0.2570075869550301
$ws.Range("E21").Value = 0.1458494226418345
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.5737522711056897
$ws.Range("H21").Value = 0.7119422934636859
$ws.Range("K21").Value = 0.3535516831559846
$ws.Range("L21").Value = 0.2043611023105285
$ws.Range("M21").Value = 0.170015111716566
$ws.Range("N21").Value = 1.671848748009047
$ws.Range("O21").Value = 2.548616017915137

$ws.Range("B22").Value = 0.6367192221275673
$ws.Range("C22").Value = 0.2562389776910727
$ws.Range("E22").Value = 0.1457149110218658
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.5712847003262524
$ws.Range("H22").Value = 0.7085103263935437
$ws.Range("K22").Value = 0.3785714037366574
$ws.Range("L22").Value = 0.2073718190695359
$ws.Range("M22").Value = 0.1761329959075155
$ws.Range("N22").Value = 1.661347594389337
$ws.Range("O22").Value = 2.536348476631531

$ws.Range("B23").Value = 0.6204750417803666
$ws.Range("C23").Value = 0.256645399867562
$ws.Range("E23").Value = 0.1457819458483272
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.5725749740439952
$ws.Range("H23").Value = 0.7103212966310153
$ws.Range("K23").Value = 0.3652228284928469
$ws.Range("L23").Value = 0.2057580611739951
$ws.Range("M23").Value = 0.1728640617341313
$ws.Range("N23").Value = 1.666910652615933
$ws.Range("O23").Value = 2.54279642216872

$ws.Range("B24").Value = 0.5590988983668126
$ws.Range("C24").Value = 0.2582682813062824
$ws.Range("E24").Value = 0.1461404729946949
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.5780484906742132
$ws.Range("H24").Value = 0.717635874044376
$ws.Range("K24").Value = 0.3145621213910204
$ws.Range("L24").Value = 0.1997978864059462
$ws.Range("M24").Value = 0.1605657063256061
$ws.Range("N24").Value = 1.688894092621982
$ws.Range("O24").Value = 2.569402658675898

$ws.Range("B25").Value = 0.493232572134616
$ws.Range("C25").Value = 0.2601985007193583
$ws.Range("E25").Value = 0.1467516956276533
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.5852056004821335
$ws.Range("H25").Value = 0.7265014578613105
$ws.Range("K25").Value = 0.2597108832355843
$ws.Range("L25").Value = 0.1957075893008167
$ws.Range("M25").Value = 0.1518641099739924
$ws.Range("N25").Value = 1.714561443615949
$ws.Range("O25").Value = 2.590907838312006
